$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows right after row 34 ("ISG expressing immune cells"),
# pushing the rest of the table down by 2 rows. Formatting is copied from
# the row above (row 34), matching Excel's default insert behavior.
$ws.Rows.Item(35).Resize(2).Insert()

# Populate the two new rows. Column order mirrors the authored edit:
# geneSymbolmore1 (C) for both rows, then cellName (B) for both rows,
# then shortName (E) for both rows, with tissueType (A) reusing the
# existing "Immune system" shared string.
$ws.Range("C35").Value = "CD3,CD4,CD69,CD154"
$ws.Range("C36").Value = "CD3,CD8,CD69,CD137"
$ws.Range("B35").Value = "Activated CD4+ T cells"
$ws.Range("B36").Value = "Activated CD8+ T cells"
$ws.Range("E35").Value = "Activated CD4+"
$ws.Range("E36").Value = "Activated CD8+"
$ws.Range("A35").Value = "Immune system"
$ws.Range("A36").Value = "Immune system"

# Update the saved selection / view state to match the post-edit session
# (entire row 57 selected, as when clicking a row header).
$ws.Rows.Item(57).Select()
